$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 58
$ws.Range("H58").Value2 = 488.85715
$ws.Range("I58").Value2 = 488.85715
$ws.Range("K58").Value2 = 1466.57145
$ws.Range("M58").Value2 = -1316.57145
# row 98
$ws.Range("H98").Value2 = 4838.8696
$ws.Range("I98").Value2 = 4466.3
$ws.Range("K98").Value2 = 4466.3
$ws.Range("M98").Value2 = -2968.3
# row 103
$ws.Range("H103").Value2 = 4763281
$ws.Range("I103").Value2 = 2261.5
$ws.Range("J103").Value2 = 10204446
$ws.Range("K103").Value2 = 6784.5
$ws.Range("L103").Value2 = 30613338
$ws.Range("M103").Value2 = -6198.5
$ws.Range("N103").Value2 = -30614510
# row 112
$ws.Range("H112").Value2 = 1462.8572
$ws.Range("J112").Value2 = 2185
$ws.Range("L112").Value2 = 6555
$ws.Range("N112").Value2 = -8771
# row 122
$ws.Range("H122").Value2 = 4838.8696
$ws.Range("I122").Value2 = 4466.3
$ws.Range("K122").Value2 = 13398.9
$ws.Range("M122").Value2 = -10948.9
# row 137
$ws.Range("H137").Value2 = 7735.174
$ws.Range("I137").Value2 = 2540.6365
$ws.Range("K137").Value2 = 7621.9095
$ws.Range("M137").Value2 = -5071.9095
# row 138
$ws.Range("H138").Value2 = 2544.6567
$ws.Range("I138").Value2 = 1500
$ws.Range("J138").Value2 = 2560.4849
$ws.Range("K138").Value2 = 4500
$ws.Range("L138").Value2 = 7681.4547
$ws.Range("M138").Value2 = 640
$ws.Range("N138").Value2 = -17961.4547

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 61
$ws.Range("H61").Value2 = 3455.25
$ws.Range("I61").Value2 = 1862.25
$ws.Range("J61").Value2 = 4251.75
$ws.Range("K61").Value2 = 1862.25
$ws.Range("L61").Value2 = 4251.75
$ws.Range("M61").Value2 = -1650.25
$ws.Range("N61").Value2 = -4675.75
# row 74
$ws.Range("H74").Value2 = 243798.48
$ws.Range("I74").Value2 = 266741.2
$ws.Range("K74").Value2 = 266741.2
$ws.Range("M74").Value2 = -265867.2
# row 77
$ws.Range("H77").Value2 = 243798.48
$ws.Range("I77").Value2 = 266741.2
$ws.Range("K77").Value2 = 1333706
$ws.Range("M77").Value2 = -1329338
# row 113
$ws.Range("H113").Value2 = 89999
$ws.Range("I113").Value2 = 0
$ws.Range("J113").Value2 = 89999
$ws.Range("K113").Value2 = 0
$ws.Range("L113").Value2 = 89999
$ws.Range("M113").Value2 = $null
$ws.Range("N113").Value2 = -98677
# row 132
$ws.Range("H132").Value2 = 1825.3214
$ws.Range("I132").Value2 = 1121.579
$ws.Range("K132").Value2 = 3364.737
$ws.Range("M132").Value2 = -834.7370000000001
# row 136
$ws.Range("H136").Value2 = 3455.25
$ws.Range("I136").Value2 = 1862.25
$ws.Range("J136").Value2 = 4251.75
$ws.Range("K136").Value2 = 5586.75
$ws.Range("L136").Value2 = 12755.25
$ws.Range("M136").Value2 = -3036.75
$ws.Range("N136").Value2 = -17855.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 20
$ws.Range("H20").Value2 = 50007976
$ws.Range("I20").Value2 = 62509224
$ws.Range("K20").Value2 = 62509224
$ws.Range("M20").Value2 = -62508977
# row 64
$ws.Range("H64").Value2 = 1477
$ws.Range("I64").Value2 = 1005.5
$ws.Range("K64").Value2 = 1005.5
$ws.Range("M64").Value2 = -780.5
# row 67
$ws.Range("H67").Value2 = 1477
$ws.Range("I67").Value2 = 1005.5
$ws.Range("K67").Value2 = 1005.5
$ws.Range("M67").Value2 = -225.5
# row 86
$ws.Range("H86").Value2 = 2055.0605
$ws.Range("I86").Value2 = 2000.8462
$ws.Range("J86").Value2 = 2256.4285
$ws.Range("K86").Value2 = 2000.8462
$ws.Range("L86").Value2 = 2256.4285
$ws.Range("M86").Value2 = -877.8462
$ws.Range("N86").Value2 = -4502.4285
# row 89
$ws.Range("H89").Value2 = 2055.0605
$ws.Range("I89").Value2 = 2000.8462
$ws.Range("J89").Value2 = 2256.4285
$ws.Range("K89").Value2 = 10004.231
$ws.Range("L89").Value2 = 11282.1425
$ws.Range("M89").Value2 = -4388.231
$ws.Range("N89").Value2 = -22514.1425
# row 94
$ws.Range("H94").Value2 = 117654010
$ws.Range("J94").Value2 = 1477.7142
$ws.Range("L94").Value2 = 1477.7142
$ws.Range("N94").Value2 = -2379.7142
# row 99
$ws.Range("H99").Value2 = 65584.69
$ws.Range("I99").Value2 = 74275.42999999999
$ws.Range("K99").Value2 = 74275.42999999999
$ws.Range("M99").Value2 = -72777.42999999999
# row 134
$ws.Range("H134").Value2 = 1878.4517
$ws.Range("I134").Value2 = 1229.9166
$ws.Range("K134").Value2 = 3689.7498
$ws.Range("M134").Value2 = -1154.7498

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value2 = 3861.7954
$ws.Range("I31").Value2 = 2715.7693
$ws.Range("J31").Value2 = 5517.1665
$ws.Range("K31").Value2 = 2715.7693
$ws.Range("L31").Value2 = 5517.1665
$ws.Range("M31").Value2 = -2420.7693
$ws.Range("N31").Value2 = -6107.1665
# row 34
$ws.Range("H34").Value2 = 3861.7954
$ws.Range("I34").Value2 = 2715.7693
$ws.Range("J34").Value2 = 5517.1665
$ws.Range("K34").Value2 = 2715.7693
$ws.Range("L34").Value2 = 5517.1665
$ws.Range("M34").Value2 = -2513.7693
$ws.Range("N34").Value2 = -5921.1665
# row 58
$ws.Range("H58").Value2 = 2192.5557
$ws.Range("I58").Value2 = 1573.3846
$ws.Range("K58").Value2 = 1573.3846
$ws.Range("M58").Value2 = -1370.3846
# row 70
$ws.Range("H70").Value2 = 49999
$ws.Range("I70").Value2 = 0
$ws.Range("K70").Value2 = 0
$ws.Range("M70").Value2 = $null
# row 73
$ws.Range("H73").Value2 = 49999
$ws.Range("I73").Value2 = 0
$ws.Range("K73").Value2 = 0
$ws.Range("M73").Value2 = $null
# row 97
$ws.Range("H97").Value2 = 44994
$ws.Range("J97").Value2 = 44994
$ws.Range("L97").Value2 = 44994
$ws.Range("N97").Value2 = -46976
# row 107
$ws.Range("H107").Value2 = 5000989
$ws.Range("I107").Value2 = 7143498.5
$ws.Range("J107").Value2 = 1799
$ws.Range("K107").Value2 = 7143498.5
$ws.Range("L107").Value2 = 1799
$ws.Range("M107").Value2 = -7141578.5
$ws.Range("N107").Value2 = -5639
# row 109
$ws.Range("H109").Value2 = 0
$ws.Range("J109").Value2 = 0
$ws.Range("L109").Value2 = 0
$ws.Range("N109").Value2 = $null
# row 112
$ws.Range("H112").Value2 = 0
$ws.Range("J112").Value2 = 0
$ws.Range("L112").Value2 = 0
$ws.Range("N112").Value2 = $null
# row 115
$ws.Range("H115").Value2 = 0
$ws.Range("J115").Value2 = 0
$ws.Range("L115").Value2 = 0
$ws.Range("N115").Value2 = $null
# row 132
$ws.Range("H132").Value2 = 17549856
$ws.Range("I132").Value2 = 5433.5557
$ws.Range("J132").Value2 = 33339834
$ws.Range("K132").Value2 = 16300.6671
$ws.Range("L132").Value2 = 100019502
$ws.Range("M132").Value2 = -13770.6671
$ws.Range("N132").Value2 = -100024562
# row 136
$ws.Range("H136").Value2 = 2192.5557
$ws.Range("I136").Value2 = 1573.3846
$ws.Range("K136").Value2 = 4720.1538
$ws.Range("M136").Value2 = -2170.1538

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 12
$ws.Range("H12").Value2 = 3994
$ws.Range("J12").Value2 = 3994
$ws.Range("L12").Value2 = 11982
$ws.Range("N12").Value2 = -12328

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value2 = 5963.5557
$ws.Range("I70").Value2 = 4000
$ws.Range("J70").Value2 = 6209
$ws.Range("K70").Value2 = 4000
$ws.Range("L70").Value2 = 6209
$ws.Range("M70").Value2 = -3730
$ws.Range("N70").Value2 = -6749
# row 73
$ws.Range("H73").Value2 = 5963.5557
$ws.Range("I73").Value2 = 4000
$ws.Range("J73").Value2 = 6209
$ws.Range("K73").Value2 = 4000
$ws.Range("L73").Value2 = 6209
$ws.Range("M73").Value2 = -3064
$ws.Range("N73").Value2 = -8081
# row 97
$ws.Range("H97").Value2 = 1337.8889
$ws.Range("I97").Value2 = 1460.875
$ws.Range("K97").Value2 = 1460.875
$ws.Range("M97").Value2 = -964.875
# row 122
$ws.Range("H122").Value2 = 3500265.2
$ws.Range("I122").Value2 = 5497707.5
$ws.Range("J122").Value2 = 4741.375
$ws.Range("K122").Value2 = 16493122.5
$ws.Range("L122").Value2 = 14224.125
$ws.Range("M122").Value2 = -16490672.5
$ws.Range("N122").Value2 = -19124.125
# row 132
$ws.Range("H132").Value2 = 2974.0588
$ws.Range("I132").Value2 = 2858
$ws.Range("J132").Value2 = 3139.8572
$ws.Range("K132").Value2 = 8574
$ws.Range("L132").Value2 = 9419.571599999999
$ws.Range("M132").Value2 = -6044
$ws.Range("N132").Value2 = -14479.5716

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 55
$ws.Range("H55").Value2 = 554.2857
$ws.Range("I55").Value2 = 566.75
$ws.Range("K55").Value2 = 566.75
$ws.Range("M55").Value2 = -393.75
# row 93
$ws.Range("H93").Value2 = 395
$ws.Range("J93").Value2 = 0
$ws.Range("L93").Value2 = 0
$ws.Range("N93").Value2 = $null
# row 128
$ws.Range("H128").Value2 = 0
$ws.Range("J128").Value2 = 0
$ws.Range("L128").Value2 = 0
$ws.Range("N128").Value2 = $null
# row 132
$ws.Range("H132").Value2 = 7513.524
$ws.Range("I132").Value2 = 5995.385
$ws.Range("J132").Value2 = 9980.5
$ws.Range("K132").Value2 = 17986.155
$ws.Range("L132").Value2 = 29941.5
$ws.Range("M132").Value2 = -15456.155
$ws.Range("N132").Value2 = -35001.5
# row 136
$ws.Range("H136").Value2 = 4758.5
$ws.Range("I136").Value2 = 3494.2
$ws.Range("J136").Value2 = 6022.8
$ws.Range("K136").Value2 = 10482.6
$ws.Range("L136").Value2 = 18068.4
$ws.Range("M136").Value2 = -7932.599999999999
$ws.Range("N136").Value2 = -23168.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 43
$ws.Range("H43").Value2 = 39980
$ws.Range("I43").Value2 = 0
$ws.Range("J43").Value2 = 39980
$ws.Range("K43").Value2 = 0
$ws.Range("L43").Value2 = 39980
$ws.Range("M43").Value2 = $null
$ws.Range("N43").Value2 = -40278
# row 122
$ws.Range("H122").Value2 = 16671155
$ws.Range("I122").Value2 = 4963.231
$ws.Range("K122").Value2 = 14889.693
$ws.Range("M122").Value2 = -12439.693
# row 132
$ws.Range("H132").Value2 = 1330.5834
$ws.Range("I132").Value2 = 1199.289
$ws.Range("J132").Value2 = 3300
$ws.Range("K132").Value2 = 3597.867
$ws.Range("L132").Value2 = 9900
$ws.Range("M132").Value2 = -1067.867
$ws.Range("N132").Value2 = -14960
